# Apputil class fix: correct the "Delete Customer" result sheet.
#
# The previous test run left a stale "Customer-00000010026" row (highlighted
# with a gray/white "selected" style) at the top of the customer list. The
# corrected result file drops that stale row, so the remaining customer
# numbers move up one slot and pick up the newer customer numbers produced
# by the latest run; the old "selected" highlight goes away too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Delete Customer")

# Row that used to be highlighted (Customer-00000010026) now shows what used
# to be the second customer, and rows 3/4 pick up the refreshed numbers.
$ws.Range("A2").Value = "Customer-00000006849"
$ws.Range("A3").Value = "Customer-00000011271"
$ws.Range("A4").Value = "Customer-00000011286"

# A2 no longer represents a "selected" row, so drop the gray-fill/white-font
# highlight and restore the plain header-style formatting (same as A1/B1).
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
